$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.199.20"
$ws.Range("E2").Value = "  +3.06%  "
$ws.Range("D3").Value = "1.581.53"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'212.14"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "'0.510"
$ws.Range("E6").Value = "  +6.13%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'26.34"
$ws.Range("E8").Value = "  +10.28%  "
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").Value = "'0.0905"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "1.807.01"
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").Value = "1.569.58"
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").Value = "29.232.26"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("D15").Value = "'3.72"
$ws.Range("E15").Value = "  +2.89%  "
$ws.Range("D16").Value = "'0.524"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("D17").Value = "'62.93"
$ws.Range("E17").Value = "  +3.37%  "
$ws.Range("D18").Value = "'238.08"
$ws.Range("E18").Value = "  +4.60%  "
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "'3.99"
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("E23").Value = "  +3.25%  "
$ws.Range("E24").Value = "  +2.96%  "
$ws.Range("D25").Value = "'154.54"
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("E26").Value = "  +4.92%  "
$ws.Range("E27").Value = "  +2.65%  "
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("D33").Value = "1.421.89"
$ws.Range("E33").Value = "  +2.34%  "
$ws.Range("E34").Value = "  +2.08%  "
$ws.Range("E35").Value = "  -3.38%  "
$ws.Range("D36").Value = "'2.83"
$ws.Range("E36").Value = "  +9.59%  "
$ws.Range("E37").Value = "  +1.36%  "
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("E39").Value = "  +2.00%  "
$ws.Range("D40").Value = "'0.530"
$ws.Range("E40").Value = "  +3.50%  "
$ws.Range("E41").Value = "  +2.39%  "
$ws.Range("D42").Value = "'53.49"
$ws.Range("E42").Value = "  +25.66%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.796"
$ws.Range("E43").Value = "  +2.63%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("D45").Value = "'0.0470"
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("D46").Value = "'64.49"
$ws.Range("E46").Value = "  +4.12%  "
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("D48").Value = "1.718.80"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").Value = "'0.839"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("D50").Value = "'85.69"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0517"
$ws.Range("E51").Value = "  +1.41%  "
